$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-10 from 45212 to 45221
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
